$d = $word.ActiveDocument

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2, wdHeaderFooterEvenPages = 3
$hfTypes = @(1, 2, 3)

# Swap the two logo picture "file names" that got mixed up between the
# default/primary and first-page headers & footers. The pictures are
# identified by their (stable) description text, since that is what
# reliably round-trips through InlineShape.AlternativeText:
#   BTec_Logo-Orange picture -> renamed to image2.jpg
#   PearsonLogo picture      -> renamed to image1.png
$renameByDescription = @{
    "BTec_Logo-Orange" = "image2.jpg"
    "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" = "image1.png"
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    foreach ($t in $hfTypes) {
        $hdr = $sec.Headers.Item($t)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                $newName = $renameByDescription[$shp.AlternativeText]
                if ($newName) {
                    $shp.Name = $newName
                }
            }
        }

        $ftr = $sec.Footers.Item($t)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                $newName = $renameByDescription[$shp.AlternativeText]
                if ($newName) {
                    $shp.Name = $newName
                }
            }
        }
    }
}
